$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New xG_home / xG_away / goals_home / goals_away values for the matches
# that were previously missing them (rows 10-15, matching Bologna's
# remaining fixtures). Values are written as text (shared-string) cells,
# matching the rest of the sheet's convention of storing these numeric
# looking values as strings rather than numbers.

function Set-TextValue($cell, $value) {
    $rng = $ws.Range($cell)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

$rows = @{
    10 = @{ D = "1.80058";  E = "0.401219"; F = "1"; G = "0" }
    11 = @{ D = "1.98287";  E = "1.1214";   F = "3"; G = "1" }
    12 = @{ D = "0.354826"; E = "4.03543";  F = "1"; G = "5" }
    13 = @{ D = "2.66532";  E = "1.61165";  F = "2"; G = "2" }
    14 = @{ D = "1.10564";  E = "2.17731";  F = "1"; G = "1" }
    15 = @{ D = "1.08799";  E = "1.89381";  F = "2"; G = "2" }
}

foreach ($r in $rows.Keys) {
    $cols = $rows[$r]
    foreach ($c in $cols.Keys) {
        Set-TextValue "$c$r" $cols[$c]
    }
}

Write-Output "done"
